$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1764705882352941
$ws.Range("C2").Value = 0.6055363321799307
$ws.Range("P2").Value = 0.1176470588235294
$ws.Range("S2").Value = 0.1003460207612457
$ws.Range("B3").Value = 0.01630434782608696
$ws.Range("C3").Value = 0.03804347826086957
$ws.Range("P3").Value = 0.7663043478260869
$ws.Range("S3").Value = 0.1793478260869565
$ws.Range("J4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.7346938775510204
$ws.Range("S4").Value = 0.2448979591836735
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.1012658227848101
$ws.Range("D6").Value = 0.02109704641350211
$ws.Range("F6").Value = 0.0379746835443038
$ws.Range("J6").Value = 0.2151898734177215
$ws.Range("O6").Value = 0.01687763713080169
$ws.Range("Q6").Value = 0.1645569620253164
$ws.Range("R6").Value = 0.1181434599156118
$ws.Range("S6").Value = 0.3248945147679325
$ws.Range("B7").Value = 0.1242603550295858
$ws.Range("D7").Value = 0.02958579881656805
$ws.Range("F7").Value = 0.07100591715976332
$ws.Range("J7").Value = 0.106508875739645
$ws.Range("O7").Value = 0.01775147928994083
$ws.Range("Q7").Value = 0.1715976331360947
$ws.Range("R7").Value = 0.08875739644970414
$ws.Range("S7").Value = 0.3905325443786982
$ws.Range("B8").Value = 0.07708333333333334
$ws.Range("D8").Value = 0.02083333333333333
$ws.Range("E8").Value = 0.002083333333333333
$ws.Range("F8").Value = 0.07083333333333333
$ws.Range("J8").Value = 0.1229166666666667
$ws.Range("O8").Value = 0.02083333333333333
$ws.Range("Q8").Value = 0.1729166666666667
$ws.Range("R8").Value = 0.09583333333333334
$ws.Range("S8").Value = 0.4166666666666667
$ws.Range("B9").Value = 0.08016877637130802
$ws.Range("D9").Value = 0.02953586497890295
$ws.Range("E9").Value = 0.004219409282700422
$ws.Range("F9").Value = 0.06329113924050633
$ws.Range("J9").Value = 0.1561181434599156
$ws.Range("O9").Value = 0.01265822784810127
$ws.Range("Q9").Value = 0.1772151898734177
$ws.Range("R9").Value = 0.05063291139240506
$ws.Range("S9").Value = 0.4261603375527426
$ws.Range("B10").Value = 0.0978013646702047
$ws.Range("D10").Value = 0.01895375284306293
$ws.Range("E10").Value = 0.003032600454890068
$ws.Range("F10").Value = 0.08112206216830932
$ws.Range("J10").Value = 0.1182714177407127
$ws.Range("O10").Value = 0.01288855193328279
$ws.Range("Q10").Value = 0.2145564821834723
$ws.Range("R10").Value = 0.1053828658074299
$ws.Range("S10").Value = 0.3479909021986353
$ws.Range("G11").Value = 0.1372549019607843
$ws.Range("J11").Value = 0.08627450980392157
$ws.Range("K11").Value = 0.1686274509803922
$ws.Range("L11").Value = 0.596078431372549
$ws.Range("S11").Value = 0.01176470588235294
$ws.Range("G12").Value = 0.7341772151898734
$ws.Range("J12").Value = 0.1835443037974684
$ws.Range("L12").Value = 0.04430379746835443
$ws.Range("S12").Value = 0.0379746835443038
$ws.Range("G13").Value = 0.6285714285714286
$ws.Range("J13").Value = 0.3142857142857143
$ws.Range("S13").Value = 0.05714285714285714
$ws.Range("F15").Value = 0.01310043668122271
$ws.Range("H15").Value = 0.1572052401746725
$ws.Range("I15").Value = 0.07423580786026202
$ws.Range("J15").Value = 0.3973799126637554
$ws.Range("K15").Value = 0.06986899563318777
$ws.Range("N15").Value = 0.004366812227074236
$ws.Range("O15").Value = 0.07423580786026202
$ws.Range("S15").Value = 0.2096069868995633
$ws.Range("F16").Value = 0.01485148514851485
$ws.Range("H16").Value = 0.1930693069306931
$ws.Range("I16").Value = 0.1089108910891089
$ws.Range("J16").Value = 0.3960396039603961
$ws.Range("K16").Value = 0.08415841584158416
$ws.Range("M16").Value = 0.009900990099009901
$ws.Range("O16").Value = 0.0396039603960396
$ws.Range("S16").Value = 0.1534653465346535
$ws.Range("F17").Value = 0.0148619957537155
$ws.Range("H17").Value = 0.1656050955414013
$ws.Range("I17").Value = 0.118895966029724
$ws.Range("J17").Value = 0.4607218683651805
$ws.Range("K17").Value = 0.07218683651804671
$ws.Range("M17").Value = 0.0148619957537155
$ws.Range("N17").Value = 0.002123142250530786
$ws.Range("O17").Value = 0.04883227176220807
$ws.Range("S17").Value = 0.1019108280254777
$ws.Range("F18").Value = 0.01680672268907563
$ws.Range("H18").Value = 0.180672268907563
$ws.Range("I18").Value = 0.1008403361344538
$ws.Range("J18").Value = 0.361344537815126
$ws.Range("K18").Value = 0.09243697478991597
$ws.Range("M18").Value = 0.01260504201680672
$ws.Range("O18").Value = 0.07563025210084033
$ws.Range("S18").Value = 0.1596638655462185
$ws.Range("F19").Value = 0.01251956181533646
$ws.Range("H19").Value = 0.2237871674491393
$ws.Range("I19").Value = 0.09389671361502347
$ws.Range("J19").Value = 0.3732394366197183
$ws.Range("K19").Value = 0.09467918622848201
$ws.Range("M19").Value = 0.01956181533646323
$ws.Range("N19").Value = 0.000782472613458529
$ws.Range("O19").Value = 0.0782472613458529
$ws.Range("S19").Value = 0.1032863849765258

Write-Host "Applied team matrix updates for games pulled March 7"
